$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 44.30160133333333
$ws.Range("H2").Value = 132.904804
$ws.Range("I2").Value = 0.1310981746002036
$ws.Range("J2").Value = 0.1395903267786693
$ws.Range("M2").Value = 5.188164
$ws.Range("N2").Value = 15.564492
$ws.Range("O2").Value = 0.133836950455521
$ws.Range("P2").Value = 0.1378024441825222
$ws.Range("Q2").Value = 229.843973179952
$ws.Range("R2").Value = 2068.595758619568
$ws.Range("S2").Value = 0.01754577989877669
$ws.Range("T2").Value = 0.01923588821433762

# Row 3
$ws.Range("G3").Value = 44.30160133333333
$ws.Range("H3").Value = 132.904804
$ws.Range("I3").Value = 0.1310981746002036
$ws.Range("J3").Value = 0.1395903267786693
$ws.Range("O3").Value = 0.1454613076012213
$ws.Range("P3").Value = 0.149771222769273
$ws.Range("Q3").Value = 249.8069835663733
$ws.Range("R3").Value = 2248.262852097359
$ws.Range("S3").Value = 0.01906971190147883
$ws.Range("T3").Value = 0.0209066139284037

# Row 4
$ws.Range("G4").Value = 44.30160133333333
$ws.Range("H4").Value = 132.904804
$ws.Range("I4").Value = 0.1310981746002036
$ws.Range("J4").Value = 0.1395903267786693
$ws.Range("M4").Value = 10.38992866666667
$ws.Range("N4").Value = 31.169786
$ws.Range("O4").Value = 0.268024751761329
$ws.Range("P4").Value = 0.2759661346766835
$ws.Range("Q4").Value = 460.2904776724383
$ws.Range("R4").Value = 4142.614299051944
$ws.Range("S4").Value = 0.03513755570358294
$ws.Range("T4").Value = 0.03852220291936453

# Row 5
$ws.Range("G5").Value = 44.30160133333333
$ws.Range("H5").Value = 132.904804
$ws.Range("I5").Value = 0.1310981746002036
$ws.Range("J5").Value = 0.1395903267786693
$ws.Range("M5").Value = 3.346565
$ws.Range("N5").Value = 6.69313
$ws.Range("O5").Value = 0.08632997224088917
$ws.Range("P5").Value = 0.05925857864370806
$ws.Range("Q5").Value = 148.2581884660867
$ws.Range("R5").Value = 889.5491307965199
$ws.Range("S5").Value = 0.01131770177406682
$ws.Range("T5").Value = 0.008271924357314685

# Row 6
$ws.Range("G6").Value = 44.30160133333333
$ws.Range("H6").Value = 132.904804
$ws.Range("I6").Value = 0.1310981746002036
$ws.Range("J6").Value = 0.1395903267786693
$ws.Range("M6").Value = 14.20137266666667
$ws.Range("N6").Value = 42.604118
$ws.Range("O6").Value = 0.3663470179410397
$ws.Range("P6").Value = 0.3772016197278132
$ws.Range("Q6").Value = 629.1435502647635
$ws.Range("R6").Value = 5662.291952382871
$ws.Range("S6").Value = 0.04802742532229835
$ws.Range("T6").Value = 0.05265369735924881

# Row 7
$ws.Range("I7").Value = 0.2142454163706631
$ws.Range("J7").Value = 0.2281236010586413
$ws.Range("M7").Value = 5.188164
$ws.Range("N7").Value = 15.564492
$ws.Range("O7").Value = 0.133836950455521
$ws.Range("P7").Value = 0.1378024441825222
$ws.Range("Q7").Value = 375.619400379888
$ws.Range("R7").Value = 3380.574603418992
$ws.Range("S7").Value = 0.0286739531761229
$ws.Range("T7").Value = 0.03143598980159938

# Row 8
$ws.Range("I8").Value = 0.2142454163706631
$ws.Range("J8").Value = 0.2281236010586413
$ws.Range("O8").Value = 0.1454613076012213
$ws.Range("P8").Value = 0.149771222769273
$ws.Range("S8").Value = 0.03116441841284475
$ws.Range("T8").Value = 0.03416635067308253

# Row 9
$ws.Range("I9").Value = 0.2142454163706631
$ws.Range("J9").Value = 0.2281236010586413
$ws.Range("M9").Value = 10.38992866666667
$ws.Range("N9").Value = 31.169786
$ws.Range("O9").Value = 0.268024751761329
$ws.Range("P9").Value = 0.2759661346766835
$ws.Range("Q9").Value = 752.2234793971708
$ws.Range("R9").Value = 6770.011314574536
$ws.Range("S9").Value = 0.05742307453874956
$ws.Range("T9").Value = 0.06295438841267902

# Row 10
$ws.Range("I10").Value = 0.2142454163706631
$ws.Range("J10").Value = 0.2281236010586413
$ws.Range("M10").Value = 3.346565
$ws.Range("N10").Value = 6.69313
$ws.Range("O10").Value = 0.08632997224088917
$ws.Range("P10").Value = 0.05925857864370806
$ws.Range("Q10").Value = 242.28893663198
$ws.Range("R10").Value = 1453.73361979188
$ws.Range("S10").Value = 0.01849580084801709
$ws.Range("T10").Value = 0.01351828035381938

# Row 11
$ws.Range("I11").Value = 0.2142454163706631
$ws.Range("J11").Value = 0.2281236010586413
$ws.Range("M11").Value = 14.20137266666667
$ws.Range("N11").Value = 42.604118
$ws.Range("O11").Value = 0.3663470179410397
$ws.Range("P11").Value = 0.3772016197278132
$ws.Range("Q11").Value = 1028.169326494819
$ws.Range("R11").Value = 9253.523938453369
$ws.Range("S11").Value = 0.07848816939492885
$ws.Range("T11").Value = 0.08604859181746097

# Row 12
$ws.Range("G12").Value = 82.35175066666666
$ws.Range("H12").Value = 247.055252
$ws.Range("I12").Value = 0.2436969288378267
$ws.Range("J12").Value = 0.2594828954344383
$ws.Range("M12").Value = 5.188164
$ws.Range("N12").Value = 15.564492
$ws.Range("O12").Value = 0.133836950455521
$ws.Range("P12").Value = 0.1378024441825222
$ws.Range("Q12").Value = 427.2543881457759
$ws.Range("R12").Value = 3845.289493311984
$ws.Range("S12").Value = 0.03261565379103083
$ws.Range("T12").Value = 0.03575737721442342

# Row 13
$ws.Range("G13").Value = 82.35175066666666
$ws.Range("H13").Value = 247.055252
$ws.Range("I13").Value = 0.2436969288378267
$ws.Range("J13").Value = 0.2594828954344383
$ws.Range("O13").Value = 0.1454613076012213
$ws.Range("P13").Value = 0.149771222769273
$ws.Range("Q13").Value = 464.3634046241866
$ws.Range("R13").Value = 4179.27064161768
$ws.Range("S13").Value = 0.03544847392715203
$ws.Range("T13").Value = 0.03886307053692722

# Row 14
$ws.Range("G14").Value = 82.35175066666666
$ws.Range("H14").Value = 247.055252
$ws.Range("I14").Value = 0.2436969288378267
$ws.Range("J14").Value = 0.2594828954344383
$ws.Range("M14").Value = 10.38992866666667
$ws.Range("N14").Value = 31.169786
$ws.Range("O14").Value = 0.268024751761329
$ws.Range("P14").Value = 0.2759661346766835
$ws.Range("Q14").Value = 855.6288150017858
$ws.Range("R14").Value = 7700.659335016072
$ws.Range("S14").Value = 0.06531680885675675
$ws.Range("T14").Value = 0.07160849166775597

# Row 15
$ws.Range("G15").Value = 82.35175066666666
$ws.Range("H15").Value = 247.055252
$ws.Range("I15").Value = 0.2436969288378267
$ws.Range("J15").Value = 0.2594828954344383
$ws.Range("M15").Value = 3.346565
$ws.Range("N15").Value = 6.69313
$ws.Range("O15").Value = 0.08632997224088917
$ws.Range("P15").Value = 0.05925857864370806
$ws.Range("Q15").Value = 275.5954864697933
$ws.Range("R15").Value = 1653.57291881876
$ws.Range("S15").Value = 0.02103834910175952
$ws.Range("T15").Value = 0.01537658756579873

# Row 16
$ws.Range("G16").Value = 82.35175066666666
$ws.Range("H16").Value = 247.055252
$ws.Range("I16").Value = 0.2436969288378267
$ws.Range("J16").Value = 0.2594828954344383
$ws.Range("M16").Value = 14.20137266666667
$ws.Range("N16").Value = 42.604118
$ws.Range("O16").Value = 0.3663470179410397
$ws.Range("P16").Value = 0.3772016197278132
$ws.Range("Q16").Value = 1169.507900969748
$ws.Range("R16").Value = 10525.57110872774
$ws.Range("S16").Value = 0.08927764316112756
$ws.Range("T16").Value = 0.09787736844953289

# Row 17
$ws.Range("G17").Value = 61.6746195
$ws.Range("H17").Value = 123.349239
$ws.Range("I17").Value = 0.1825087534596294
$ws.Range("J17").Value = 0.1295540872992837
$ws.Range("M17").Value = 5.188164
$ws.Range("N17").Value = 15.564492
$ws.Range("O17").Value = 0.133836950455521
$ws.Range("P17").Value = 0.1378024441825222
$ws.Range("Q17").Value = 319.9780406035979
$ws.Range("R17").Value = 1919.868243621588
$ws.Range("S17").Value = 0.02442641499447531
$ws.Range("T17").Value = 0.01785286988367715

# Row 18
$ws.Range("G18").Value = 61.6746195
$ws.Range("H18").Value = 123.349239
$ws.Range("I18").Value = 0.1825087534596294
$ws.Range("J18").Value = 0.1295540872992837
$ws.Range("O18").Value = 0.1454613076012213
$ws.Range("P18").Value = 0.149771222769273
$ws.Range("Q18").Value = 347.7696109442099
$ws.Range("R18").Value = 2086.61766566526
$ws.Range("S18").Value = 0.02654796192690661
$ws.Range("T18").Value = 0.01940347406957086

# Row 19
$ws.Range("G19").Value = 61.6746195
$ws.Range("H19").Value = 123.349239
$ws.Range("I19").Value = 0.1825087534596294
$ws.Range("J19").Value = 0.1295540872992837
$ws.Range("M19").Value = 10.38992866666667
$ws.Range("N19").Value = 31.169786
$ws.Range("O19").Value = 0.268024751761329
$ws.Range("P19").Value = 0.2759661346766835
$ws.Range("Q19").Value = 640.794897148809
$ws.Range("R19").Value = 3844.769382892854
$ws.Range("S19").Value = 0.04891686334028677
$ws.Range("T19").Value = 0.03575254070354894

# Row 20
$ws.Range("G20").Value = 61.6746195
$ws.Range("H20").Value = 123.349239
$ws.Range("I20").Value = 0.1825087534596294
$ws.Range("J20").Value = 0.1295540872992837
$ws.Range("M20").Value = 3.346565
$ws.Range("N20").Value = 6.69313
$ws.Range("O20").Value = 0.08632997224088917
$ws.Range("P20").Value = 0.05925857864370806
$ws.Range("Q20").Value = 206.3981230070175
$ws.Range("R20").Value = 825.59249202807
$ws.Range("S20").Value = 0.01575597561988909
$ws.Range("T20").Value = 0.007677191070838422

# Row 21
$ws.Range("G21").Value = 61.6746195
$ws.Range("H21").Value = 123.349239
$ws.Range("I21").Value = 0.1825087534596294
$ws.Range("J21").Value = 0.1295540872992837
$ws.Range("M21").Value = 14.20137266666667
$ws.Range("N21").Value = 42.604118
$ws.Range("O21").Value = 0.3663470179410397
$ws.Range("P21").Value = 0.3772016197278132
$ws.Range("Q21").Value = 875.8642555943669
$ws.Range("R21").Value = 5255.185533566202
$ws.Range("S21").Value = 0.06686153757807164
$ws.Range("T21").Value = 0.04886801157164832

# Row 22
$ws.Range("G22").Value = 77.19964866666666
$ws.Range("H22").Value = 231.598946
$ws.Range("I22").Value = 0.2284507267316773
$ws.Range("J22").Value = 0.2432490894289675
$ws.Range("M22").Value = 5.188164
$ws.Range("N22").Value = 15.564492
$ws.Range("O22").Value = 0.133836950455521
$ws.Range("P22").Value = 0.1378024441825222
$ws.Range("Q22").Value = 400.5244380250479
$ws.Range("R22").Value = 3604.719942225432
$ws.Range("S22").Value = 0.03057514859511525
$ws.Range("T22").Value = 0.03352031906848465

# Row 23
$ws.Range("G23").Value = 77.19964866666666
$ws.Range("H23").Value = 231.598946
$ws.Range("I23").Value = 0.2284507267316773
$ws.Range("J23").Value = 0.2432490894289675
$ws.Range("O23").Value = 0.1454613076012213
$ws.Range("P23").Value = 0.149771222769273
$ws.Range("Q23").Value = 435.3118349086266
$ws.Range("R23").Value = 3917.806514177639
$ws.Range("S23").Value = 0.03323074143283904
$ws.Range("T23").Value = 0.03643171356128871

# Row 24
$ws.Range("G24").Value = 77.19964866666666
$ws.Range("H24").Value = 231.598946
$ws.Range("I24").Value = 0.2284507267316773
$ws.Range("J24").Value = 0.2432490894289675
$ws.Range("M24").Value = 10.38992866666667
$ws.Range("N24").Value = 31.169786
$ws.Range("O24").Value = 0.268024751761329
$ws.Range("P24").Value = 0.2759661346766835
$ws.Range("Q24").Value = 802.0988427383952
$ws.Range("R24").Value = 7218.889584645556
$ws.Range("S24").Value = 0.06123044932195301
$ws.Range("T24").Value = 0.06712851097333508

# Row 25
$ws.Range("G25").Value = 77.19964866666666
$ws.Range("H25").Value = 231.598946
$ws.Range("I25").Value = 0.2284507267316773
$ws.Range("J25").Value = 0.2432490894289675
$ws.Range("M25").Value = 3.346565
$ws.Range("N25").Value = 6.69313
$ws.Range("O25").Value = 0.08632997224088917
$ws.Range("P25").Value = 0.05925857864370806
$ws.Range("Q25").Value = 258.3536422401633
$ws.Range("R25").Value = 1550.12185344098
$ws.Range("S25").Value = 0.01972214489715666
$ws.Range("T25").Value = 0.01441459529593685

# Row 26
$ws.Range("G26").Value = 77.19964866666666
$ws.Range("H26").Value = 231.598946
$ws.Range("I26").Value = 0.2284507267316773
$ws.Range("J26").Value = 0.2432490894289675
$ws.Range("M26").Value = 14.20137266666667
$ws.Range("N26").Value = 42.604118
$ws.Range("O26").Value = 0.3663470179410397
$ws.Range("P26").Value = 0.3772016197278132
$ws.Range("Q26").Value = 1096.34098045107
$ws.Range("R26").Value = 9867.068824059626
$ws.Range("S26").Value = 0.08369224248461332
$ws.Range("T26").Value = 0.09175395052992222
